$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet to match the batter it documents
$ws.Name = "Kuldip Yadav"

# Insert a new column before column A ("teamName" ...), shifting the
# existing teamName..result columns from A:L to B:M
$ws.Range("A1:A2").EntireColumn.Insert()

# Populate the newly inserted "matchNo" column
$ws.Range("A1").Value = "matchNo"
$ws.Range("A2").Value = "51st"
